# "reverted trial counts to normal"
#
# Columns L:O on Sheet1 held the "real" iti_min/iti_max/word_trial_count/
# nonword_trial_count values (headers a/b/c/d) while columns C:F had been
# temporarily pinned to 1. This restores C:F to the real values and removes
# the now-redundant helper columns L:O.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5: only word/nonword trial counts (E,F) need restoring.
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 20

$ws.Range("E3").Value = 30
$ws.Range("F3").Value = 30

$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 30

$ws.Range("E5").Value = 30
$ws.Range("F5").Value = 30

# Rows 6-9: iti_min/iti_max (C,D) and word/nonword trial counts (E,F)
# all need restoring.
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 8

$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 12
$ws.Range("F7").Value = 12

$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = 12
$ws.Range("F8").Value = 12

$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 12
$ws.Range("F9").Value = 12

# Drop the now-redundant helper columns (headers a, b, c, d).
$ws.Range("L1:O9").EntireColumn.Delete()

# Restore the saved selection/active-cell position.
$null = $ws.Range("P19").Select()
